$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the numeric-looking values so they stay stored
# as text (matching the original t="str" cell type) rather than being
# auto-converted to numbers by Excel.
$rng = $ws.Range("C2:F11")
$rng.NumberFormat = "@"

$values = @{
    2  = @("16", "7", "1", "1")
    3  = @("24", "11", "3", "1")
    4  = @("27", "8", "0", "4")
    5  = @("6", "9", "0", "0")
    6  = @("1", "4", "0", "0")
    7  = @("16", "10", "1", "1")
    8  = @("2", "3", "0", "0")
    9  = @("13", "3", "0", "2")
    10 = @("6", "4", "0", "1")
    11 = @("2", "4", "0", "0")
}

foreach ($row in $values.Keys) {
    $vals = $values[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
}
